# Auto-generated Excel COM-interop script to apply profit recalculation updates
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the Hades_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1857.6842
$ws.Range("I100").Value = 1662.7273
$ws.Range("J100").Value = 2125.75
$ws.Range("K100").Value = 1662.7273
$ws.Range("L100").Value = 2125.75
$ws.Range("M100").Value = -1121.7273
$ws.Range("N100").Value = -3207.75

$ws.Range("H132").Value = 682247.75
$ws.Range("I132").Value = 1654.6727
$ws.Range("J132").Value = 2884166.5
$ws.Range("K132").Value = 4964.0181
$ws.Range("L132").Value = 8652499.5
$ws.Range("M132").Value = -2434.0181
$ws.Range("N132").Value = -8657559.5

$ws.Range("H135").Value = 18973.104
$ws.Range("I135").Value = 22388.639
$ws.Range("J135").Value = 4379.4546
$ws.Range("K135").Value = 201497.751
$ws.Range("L135").Value = 39415.0914
$ws.Range("M135").Value = -198962.751
$ws.Range("N135").Value = -44485.0914

$ws.Range("H138").Value = 1917832.9
$ws.Range("I138").Value = 1459.2572
$ws.Range("J138").Value = 3207699.8
$ws.Range("K138").Value = 4377.7716
$ws.Range("L138").Value = 9623099.399999999
$ws.Range("M138").Value = 762.2284
$ws.Range("N138").Value = -9633379.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10501548
$ws.Range("I74").Value = 13211629
$ws.Range("J74").Value = 203240
$ws.Range("K74").Value = 13211629
$ws.Range("L74").Value = 203240
$ws.Range("M74").Value = -13210755
$ws.Range("N74").Value = -204988

$ws.Range("H77").Value = 10501548
$ws.Range("I77").Value = 13211629
$ws.Range("J77").Value = 203240
$ws.Range("K77").Value = 66058145
$ws.Range("L77").Value = 1016200
$ws.Range("M77").Value = -66053777
$ws.Range("N77").Value = -1024936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 544.625
$ws.Range("I94").Value = 354.15384
$ws.Range("J94").Value = 1370
$ws.Range("K94").Value = 354.15384
$ws.Range("L94").Value = 1370
$ws.Range("M94").Value = 96.84616
$ws.Range("N94").Value = -2272

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2910.325
$ws.Range("I31").Value = 1337.4
$ws.Range("J31").Value = 5531.8667
$ws.Range("K31").Value = 1337.4
$ws.Range("L31").Value = 5531.8667
$ws.Range("M31").Value = -1042.4
$ws.Range("N31").Value = -6121.8667

$ws.Range("H34").Value = 2910.325
$ws.Range("I34").Value = 1337.4
$ws.Range("J34").Value = 5531.8667
$ws.Range("K34").Value = 1337.4
$ws.Range("L34").Value = 5531.8667
$ws.Range("M34").Value = -1135.4
$ws.Range("N34").Value = -5935.8667

$ws.Range("H58").Value = 21278130
$ws.Range("I58").Value = 23257178
$ws.Range("J58").Value = 3375.25
$ws.Range("K58").Value = 23257178
$ws.Range("L58").Value = 3375.25
$ws.Range("M58").Value = -23256975
$ws.Range("N58").Value = -3781.25

$ws.Range("H112").Value = 40722
$ws.Range("J112").Value = 40722
$ws.Range("L112").Value = 40722
$ws.Range("N112").Value = -43676

$ws.Range("H132").Value = 18863.352
$ws.Range("I132").Value = 1275.8478
$ws.Range("J132").Value = 92411.09
$ws.Range("K132").Value = 3827.5434
$ws.Range("L132").Value = 277233.27
$ws.Range("M132").Value = -1297.5434
$ws.Range("N132").Value = -282293.27

$ws.Range("H134").Value = 19462.373
$ws.Range("I134").Value = 1119.4783
$ws.Range("J134").Value = 84368
$ws.Range("K134").Value = 3358.4349
$ws.Range("L134").Value = 253104
$ws.Range("M134").Value = -823.4349000000002
$ws.Range("N134").Value = -258174

$ws.Range("H136").Value = 21278130
$ws.Range("I136").Value = 23257178
$ws.Range("J136").Value = 3375.25
$ws.Range("K136").Value = 69771534
$ws.Range("L136").Value = 10125.75
$ws.Range("M136").Value = -69768984
$ws.Range("N136").Value = -15225.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 7750
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 30000
$ws.Range("M110").Value = 1090
$ws.Range("N110").Value = -38180

$ws.Range("H112").Value = 15154391
$ws.Range("I112").Value = 2239
$ws.Range("J112").Value = 20836448
$ws.Range("K112").Value = 6717
$ws.Range("L112").Value = 62509344
$ws.Range("M112").Value = -5609
$ws.Range("N112").Value = -62511560

$ws.Range("H113").Value = 595
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 595
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1785
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6125

$ws.Range("H115").Value = 2627.4666
$ws.Range("J115").Value = 2600.8572
$ws.Range("L115").Value = 7802.571599999999
$ws.Range("N115").Value = -10152.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 24131
$ws.Range("J74").Value = 24131
$ws.Range("L74").Value = 24131
$ws.Range("N74").Value = -26003

$ws.Range("H77").Value = 24131
$ws.Range("J77").Value = 24131
$ws.Range("L77").Value = 72393
$ws.Range("N77").Value = -81753

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1329.0646
$ws.Range("I100").Value = 1133.375
$ws.Range("K100").Value = 1133.375
$ws.Range("M100").Value = -592.375

$ws.Range("H132").Value = 93545.17999999999
$ws.Range("I132").Value = 4200
$ws.Range("J132").Value = 144599.58
$ws.Range("K132").Value = 12600
$ws.Range("L132").Value = 433798.74
$ws.Range("M132").Value = -10070
$ws.Range("N132").Value = -438858.74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 38347.5
$ws.Range("J75").Value = 38347.5
$ws.Range("L75").Value = 38347.5
$ws.Range("N75").Value = -40219.5

$ws.Range("H78").Value = 38347.5
$ws.Range("J78").Value = 38347.5
$ws.Range("L78").Value = 115042.5
$ws.Range("N78").Value = -124402.5

$ws.Range("H107").Value = 587.75
$ws.Range("I107").Value = 486.5
$ws.Range("J107").Value = 689
$ws.Range("K107").Value = 1459.5
$ws.Range("L107").Value = 2067
$ws.Range("M107").Value = 460.5
$ws.Range("N107").Value = -5907

$ws.Range("H126").Value = 1715.75
$ws.Range("I126").Value = 1603.7142
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4811.142599999999
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2341.142599999999
$ws.Range("N126").Value = -12440

$ws.Range("H132").Value = 42526.98
$ws.Range("I132").Value = 24058.023
$ws.Range("J132").Value = 201360
$ws.Range("K132").Value = 72174.069
$ws.Range("L132").Value = 604080
$ws.Range("M132").Value = -69644.069
$ws.Range("N132").Value = -609140

$ws.Range("H136").Value = 39462.887
$ws.Range("I136").Value = 23330.182
$ws.Range("J136").Value = 118333.89
$ws.Range("K136").Value = 69990.546
$ws.Range("L136").Value = 355001.67
$ws.Range("M136").Value = -67440.546
$ws.Range("N136").Value = -360101.67

Write-Host "Applied all profit recalculation updates across 28 rows in 8 sheets."
